# Adds three new potion entries ("Light Pink Potion", "Dark Pink Potion",
# "Stinger") right after the existing "Rainbow Potion" list item, and moves
# the "_GoBack" bookmark from "Rainbow Potion" onto the new
# "Dark Pink Potion" entry (matching the target diff).

$d = $word.ActiveDocument

function Find-ParaByText($text) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.TrimEnd([char]13) -eq $text) {
            return $p
        }
    }
    return $null
}

# 1. Locate the "Rainbow Potion" paragraph (currently the last potion).
$rainbow = Find-ParaByText "Rainbow Potion"

# 2. Insert three new list paragraphs after it. InsertParagraphAfter()
#    clones the originating paragraph's pPr/rPr (pStyle "Lijstalinea",
#    numPr ilvl=0/numId=1, rPr lang="en-GB"), exactly like the diff shows.
$null = $rainbow.Range.InsertParagraphAfter()

$pLightPink = $rainbow.Next()
$rLightPink = $pLightPink.Range
$rLightPink.End = $rLightPink.End - 1
$rLightPink.Text = "Light Pink Potion"

$null = $pLightPink.Range.InsertParagraphAfter()
$pDarkPink = $pLightPink.Next()
$rDarkPink = $pDarkPink.Range
$rDarkPink.End = $rDarkPink.End - 1
$rDarkPink.Text = "Dark Pink Potion"

$null = $pDarkPink.Range.InsertParagraphAfter()
$pStinger = $pDarkPink.Next()
$rStinger = $pStinger.Range
$rStinger.End = $rStinger.End - 1
$rStinger.Text = "Stinger"

# 3. Move the "_GoBack" bookmark off "Rainbow Potion" onto the end of the
#    new "Dark Pink Potion" run.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# A bookmark collapsed exactly at a paragraph's trailing boundary (right
# before the paragraph mark) gets relocated to the neighbouring paragraph
# by this runtime, so temporarily append a one-character sentinel run,
# anchor the bookmark just in front of it (now a genuine mid-paragraph
# position), and then delete the sentinel again.
$sentinel = $pDarkPink.Range
$sentinel.Start = $sentinel.End - 1
$sentinel.End = $sentinel.Start
$sentinel.InsertAfter("X")

$bkRange = $pDarkPink.Range
$bkRange.Start = $bkRange.End - 2
$bkRange.End = $bkRange.Start
$d.Bookmarks.Add("_GoBack", $bkRange)

$cleanup = $pDarkPink.Range
$cleanup.Start = $cleanup.End - 2
$cleanup.End = $cleanup.End - 1
$cleanup.Text = ""

$d.Save()
